# Weekly update: insert 4 new price rows (week of 2022-07-11, serial 44753)
# at the top of the "Rosara/Asterix/Rodeo" 1a (guarda) block, pushing the
# existing rows 618:649 down to 622:653.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 618:621 (shifts 618:649 -> 622:653)
$ws.Range("A618:R621").EntireRow.Insert()

# Common (constant) values shared by the whole "Femacal de La Calera / Papa" block
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$fecha     = 44753
$codreg    = 5
$categoriaId = 100114001
$categoria = "Papa"
$unidad    = "`$/saco 25 kilos"
$kgUnidades = 25
$clasificacion = "Hortaliza"

# New rows: Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg
$newRows = @(
    @{ Row = 618; Variedad = "Asterix"; Calidad = "1a (guarda)"; Volumen = 480; PMin = 7800; PMax = 8000; PProm = 7904; Origen = "Provincia de Talca";    PKg = 316 },
    @{ Row = 619; Variedad = "Rodeo";   Calidad = "1a (guarda)"; Volumen = 250; PMin = 7500; PMax = 7500; PProm = 7500; Origen = "Provincia de Talca";    PKg = 300 },
    @{ Row = 620; Variedad = "Rosara";  Calidad = "1a (guarda)"; Volumen = 410; PMin = 7500; PMax = 8000; PProm = 7780; Origen = "Provincia de Quillota"; PKg = 311 },
    @{ Row = 621; Variedad = "Rosara";  Calidad = "1a (guarda)"; Volumen = 480; PMin = 7500; PMax = 7800; PProm = 7662; Origen = "Provincia de Talca";    PKg = 306 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $categoriaId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}

Write-Output "Inserted 4 rows (618-621) and populated new weekly data."
